# Update the "dSF" column (F) values on Sheet1 to reflect a repull of
# data / recalculated mean (per commit message: "repull data, push all
# data, mean calculation"). Only the rows whose dSF value changed in the
# source diff are touched; everything else is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 3
    3  = 3
    4  = -3
    5  = 10
    6  = -1
    7  = 0
    8  = 2
    9  = -7
    11 = 3
    12 = -6
    14 = 3
    15 = 2
    18 = -3
    19 = -1
    20 = 0
    21 = 7
    22 = -3
    23 = 6
    24 = 1
    25 = 4
    26 = 0
    27 = 3
    29 = 2
    31 = 1
    32 = 0
    33 = -1
    34 = -3
    36 = 3
    37 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
